# Generate Report for Handback
# Update the timestamp values on the Overview, zh-cn, and de-de sheets
# to reflect the new handback report generation times.
#
# "Latest HO Xliff Generate Date" (Overview!G2 and de-de!H2) is the same
# underlying value on both sheets, so both must be updated together.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 19:03:54"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 19:03:49"
$wsZhCn.Range("K2").Value = "2016-08-17 19:04:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 19:03:54"
$wsDeDe.Range("K2").Value = "2016-08-17 19:04:20"
